$d = $word.ActiveDocument

# 1. Update the main heading/title text (appears twice in the document:
#    once near the top, once as the section heading before the body).
$d.Content.Find.Execute(
    "Inspiración Web: Portafolios Excepcionales en HTML, CSS y JS",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Sitios web increíbles recomendados: Portafolios Excepcionales en HTML, CSS y JS",
    2)

# 2. Remove three "Publicaciones Similares" list items that no longer apply:
#    "El Aborto", "El Mercantilismo", "Comandos De Blogdown".
#    Each is a whole list paragraph (two hyperlinks + a separating space run).
#    Delete from bottom-most to top-most so earlier paragraph indices remain valid.

$targets = @("El Mercantilismo", "Comandos De Blogdown", "El Aborto")

foreach ($target in $targets) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs.Item($i)
        $text = $p.Range.Text.Trim()
        if ($text -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}

Write-Output "done"
